$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: one between "ResultUnit" (V) and "MeasuredEnzyme" (old W),
# and one right after "MeasuredEnzyme" (which has shifted to X) and before
# "MolarExtinctionCoefficient" (old X, now Z).
$ws.Columns("W").Insert()
$ws.Columns("Y").Insert()

# New column W: TimeStep
$ws.Range("W1").Value = "TimeStep"
$ws.Range("W2").Value = "# Pas de temps de la mesure"
$ws.Range("W3").Value = "#integer,`n  unit:s"
$ws.Range("W4").Value = "# format: entier,  ne pas spécifier d'unité (seconde)"
$ws.Range("W5").Value = "# ex: ...."

# New column Y: UsedSubstrat
$ws.Range("Y1").Value = "UsedSubstrat"
$ws.Range("Y2").Value = "# substrat utilisé"
$ws.Range("Y3").Value = "#string"
$ws.Range("Y4").Value = "# format: texte"
$ws.Range("Y5").Value = "# ex: ...."

# Wavelength format note now clarifies the unit (nm)
$ws.Range("J4").Value = "# format: nombre entier, ne pas spécifier d'unité (nm)"
